# Apply the updated 漫展信息 numbers/statuses to the relevant sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 5695
$ws1.Range("G2").Value = 79

$ws1.Range("F3").Value = 380
$ws1.Range("G3").Value = "不可售"

$ws1.Range("G4").Value = "不可售"

$ws1.Range("F5").Value = 312
$ws1.Range("G5").Value = "不可售"

$ws1.Range("F7").Value = 67

$ws1.Range("F8").Value = 385

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Range("F2").Value = 53

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 5695
$ws4.Range("G2").Value = 79

$ws4.Range("F3").Value = 380
$ws4.Range("G3").Value = "不可售"

$ws4.Range("G4").Value = "不可售"

$ws4.Range("F5").Value = 312
$ws4.Range("G5").Value = "不可售"

$ws4.Range("F7").Value = 67

$ws4.Range("F8").Value = 53

$ws4.Range("F9").Value = 385
